$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original data (rows 1-11):
#   1: location | TotalDeathCount
#   2: High income | 2913655
#   3: Upper middle income | 2667294
#   4: Europe | 2082494
#   5: Asia | 1635654
#   6: North America | 1606198
#   7: South America | 1359982
#   8: Lower middle income | 1340728
#   9: Africa | 259042
#  10: Low income | 48001
#  11: Oceania | 30265
#
# Target data (rows 1-7):
#   1: location | TotalDeathCount
#   2: Europe | 2082494
#   3: Asia | 1635654
#   4: North America | 1606198
#   5: South America | 1359982
#   6: Africa | 259042
#   7: Oceania | 30265
#
# So rows 10 (Low income), 8 (Lower middle income), 3 (Upper middle
# income) and 2 (High income) need to be removed. Delete from the
# bottom up so earlier row numbers stay valid while deleting.
$ws.Rows("10:10").Delete()
$ws.Rows("8:8").Delete()
$ws.Rows("2:3").Delete()

$ws.Range("A9").Select()
